$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 39
$ws.Cells.Item(39, 8).Value = 2673.1667
$ws.Cells.Item(39, 9).Value = 1385
$ws.Cells.Item(39, 11).Value = 4155
$ws.Cells.Item(39, 13).Value = -3859

# Row 51
$ws.Cells.Item(51, 8).Value = 12260.385
$ws.Cells.Item(51, 10).Value = 8262.362999999999
$ws.Cells.Item(51, 12).Value = 8262.362999999999
$ws.Cells.Item(51, 14).Value = -9230.362999999999

# Row 70
$ws.Cells.Item(70, 8).Value = 2373
$ws.Cells.Item(70, 9).Value = 2974.25
$ws.Cells.Item(70, 11).Value = 8922.75
$ws.Cells.Item(70, 13).Value = -8652.75

# Row 73
$ws.Cells.Item(73, 8).Value = 2373
$ws.Cells.Item(73, 9).Value = 2974.25
$ws.Cells.Item(73, 11).Value = 8922.75
$ws.Cells.Item(73, 13).Value = -7986.75

# Row 87
$ws.Cells.Item(87, 8).Value = 39999.332
$ws.Cells.Item(87, 10).Value = 39999.332
$ws.Cells.Item(87, 12).Value = 39999.332
$ws.Cells.Item(87, 14).Value = -42495.332

# Row 90
$ws.Cells.Item(90, 8).Value = 39999.332
$ws.Cells.Item(90, 10).Value = 39999.332
$ws.Cells.Item(90, 12).Value = 119997.996
$ws.Cells.Item(90, 14).Value = -132477.996

# Row 94
$ws.Cells.Item(94, 8).Value = 5174.7
$ws.Cells.Item(94, 9).Value = 218.375
$ws.Cells.Item(94, 11).Value = 218.375
$ws.Cells.Item(94, 13).Value = 232.625

# Row 137
$ws.Cells.Item(137, 8).Value = 2569.45
$ws.Cells.Item(137, 9).Value = 2693.0625
$ws.Cells.Item(137, 10).Value = 2075
$ws.Cells.Item(137, 11).Value = 8079.1875
$ws.Cells.Item(137, 12).Value = 6225
$ws.Cells.Item(137, 13).Value = -5529.1875
$ws.Cells.Item(137, 14).Value = -11325

# Row 139
$ws.Cells.Item(139, 8).Value = 84329.89
$ws.Cells.Item(139, 10).Value = 84996.25
$ws.Cells.Item(139, 12).Value = 84996.25
$ws.Cells.Item(139, 14).Value = -95276.25

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Cells.Item(32, 8).Value = 1764447.6
$ws.Cells.Item(32, 9).Value = 1764447.6
$ws.Cells.Item(32, 11).Value = 1764447.6
$ws.Cells.Item(32, 13).Value = -1764160.6

# Row 74
$ws.Cells.Item(74, 8).Value = 49594.027
$ws.Cells.Item(74, 9).Value = 83039
$ws.Cells.Item(74, 10).Value = 5000.7334
$ws.Cells.Item(74, 11).Value = 83039
$ws.Cells.Item(74, 12).Value = 5000.7334
$ws.Cells.Item(74, 13).Value = -82165
$ws.Cells.Item(74, 14).Value = -6748.7334

# Row 77
$ws.Cells.Item(77, 8).Value = 49594.027
$ws.Cells.Item(77, 9).Value = 83039
$ws.Cells.Item(77, 10).Value = 5000.7334
$ws.Cells.Item(77, 11).Value = 415195
$ws.Cells.Item(77, 12).Value = 25003.667
$ws.Cells.Item(77, 13).Value = -410827
$ws.Cells.Item(77, 14).Value = -33739.667

# Row 102
$ws.Cells.Item(102, 8).Value = 1500
$ws.Cells.Item(102, 9).Value = 0
$ws.Cells.Item(102, 11).Value = 0
$ws.Cells.Item(102, 13).Value = $null

# Row 122
$ws.Cells.Item(122, 8).Value = 19553
$ws.Cells.Item(122, 9).Value = 26648.75
$ws.Cells.Item(122, 10).Value = 8199.799999999999
$ws.Cells.Item(122, 11).Value = 79946.25
$ws.Cells.Item(122, 12).Value = 24599.4
$ws.Cells.Item(122, 13).Value = -77496.25
$ws.Cells.Item(122, 14).Value = -29499.4

# Row 132
$ws.Cells.Item(132, 8).Value = 6021.6514
$ws.Cells.Item(132, 9).Value = 4549.346
$ws.Cells.Item(132, 10).Value = 8273.412
$ws.Cells.Item(132, 11).Value = 13648.038
$ws.Cells.Item(132, 12).Value = 24820.236
$ws.Cells.Item(132, 13).Value = -11118.038
$ws.Cells.Item(132, 14).Value = -29880.236

$ws = $wb.Worksheets.Item("BSM")
# Row 99
$ws.Cells.Item(99, 8).Value = 4788996.5
$ws.Cells.Item(99, 9).Value = 4490.9375
$ws.Cells.Item(99, 11).Value = 4490.9375
$ws.Cells.Item(99, 13).Value = -2992.9375

# Row 134
$ws.Cells.Item(134, 8).Value = 5761.7207
$ws.Cells.Item(134, 9).Value = 1768.1578
$ws.Cells.Item(134, 10).Value = 8923.291999999999
$ws.Cells.Item(134, 11).Value = 5304.4734
$ws.Cells.Item(134, 12).Value = 26769.876
$ws.Cells.Item(134, 13).Value = -2769.4734
$ws.Cells.Item(134, 14).Value = -31839.876

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Cells.Item(31, 8).Value = 8554590
$ws.Cells.Item(31, 9).Value = 2581.4285
$ws.Cells.Item(31, 11).Value = 2581.4285
$ws.Cells.Item(31, 13).Value = -2286.4285

# Row 34
$ws.Cells.Item(34, 8).Value = 8554590
$ws.Cells.Item(34, 9).Value = 2581.4285
$ws.Cells.Item(34, 11).Value = 2581.4285
$ws.Cells.Item(34, 13).Value = -2379.4285

# Row 105
$ws.Cells.Item(105, 8).Value = 6497955
$ws.Cells.Item(105, 10).Value = 7382.3335
$ws.Cells.Item(105, 12).Value = 7382.3335
$ws.Cells.Item(105, 14).Value = -10876.3335

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Cells.Item(5, 8).Value = 1424.1875
$ws.Cells.Item(5, 10).Value = 1780
$ws.Cells.Item(5, 12).Value = 5340
$ws.Cells.Item(5, 14).Value = -5564

# Row 38
$ws.Cells.Item(38, 8).Value = 70.75
$ws.Cells.Item(38, 10).Value = 78.5
$ws.Cells.Item(38, 12).Value = 235.5
$ws.Cells.Item(38, 14).Value = -929.5

# Row 56
$ws.Cells.Item(56, 8).Value = 7204.2856
$ws.Cells.Item(56, 9).Value = 7204.2856
$ws.Cells.Item(56, 11).Value = 7204.2856
$ws.Cells.Item(56, 13).Value = -6674.2856

# Row 76
$ws.Cells.Item(76, 8).Value = 500001500
$ws.Cells.Item(76, 9).Value = 500001500
$ws.Cells.Item(76, 11).Value = 1500004500
$ws.Cells.Item(76, 13).Value = -1500004117

# Row 79
$ws.Cells.Item(79, 8).Value = 500001500
$ws.Cells.Item(79, 9).Value = 500001500
$ws.Cells.Item(79, 11).Value = 1500004500
$ws.Cells.Item(79, 13).Value = -1500003174

# Row 113
$ws.Cells.Item(113, 8).Value = 1518.3334
$ws.Cells.Item(113, 9).Value = 1503.25
$ws.Cells.Item(113, 10).Value = 1533.4166
$ws.Cells.Item(113, 11).Value = 4509.75
$ws.Cells.Item(113, 12).Value = 4600.2498
$ws.Cells.Item(113, 13).Value = -2339.75
$ws.Cells.Item(113, 14).Value = -8940.2498

# Row 115
$ws.Cells.Item(115, 8).Value = 1670.7273
$ws.Cells.Item(115, 9).Value = 792.6667
$ws.Cells.Item(115, 11).Value = 2378.0001
$ws.Cells.Item(115, 13).Value = -1203.0001

# Row 135
$ws.Cells.Item(135, 8).Value = 1424.1875
$ws.Cells.Item(135, 10).Value = 1780
$ws.Cells.Item(135, 12).Value = 16020
$ws.Cells.Item(135, 14).Value = -21090

# Row 137
$ws.Cells.Item(137, 8).Value = 113130.836
$ws.Cells.Item(137, 9).Value = 92399.73
$ws.Cells.Item(137, 10).Value = 145708.28
$ws.Cells.Item(137, 11).Value = 277199.19
$ws.Cells.Item(137, 12).Value = 437124.84
$ws.Cells.Item(137, 13).Value = -272099.19
$ws.Cells.Item(137, 14).Value = -447324.84

$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Cells.Item(132, 8).Value = 1756.1034
$ws.Cells.Item(132, 9).Value = 1553.9584
$ws.Cells.Item(132, 10).Value = 2726.4
$ws.Cells.Item(132, 11).Value = 4661.8752
$ws.Cells.Item(132, 12).Value = 8179.200000000001
$ws.Cells.Item(132, 13).Value = -2131.8752
$ws.Cells.Item(132, 14).Value = -13239.2

$ws = $wb.Worksheets.Item("LTW")
# Row 2
$ws.Cells.Item(2, 8).Value = 0
$ws.Cells.Item(2, 9).Value = 0
$ws.Cells.Item(2, 11).Value = 0

# Row 7
$ws.Cells.Item(7, 8).Value = 5184.5938
$ws.Cells.Item(7, 9).Value = 3995.9565
$ws.Cells.Item(7, 10).Value = 8222.223
$ws.Cells.Item(7, 11).Value = 3995.9565
$ws.Cells.Item(7, 12).Value = 8222.223
$ws.Cells.Item(7, 13).Value = -3883.9565
$ws.Cells.Item(7, 14).Value = -8446.223

# Row 93
$ws.Cells.Item(93, 8).Value = 5595
$ws.Cells.Item(93, 10).Value = 12726
$ws.Cells.Item(93, 12).Value = 12726
$ws.Cells.Item(93, 14).Value = -15222

# Row 122
$ws.Cells.Item(122, 8).Value = 4421
$ws.Cells.Item(122, 9).Value = 2872.3333
$ws.Cells.Item(122, 11).Value = 8616.999899999999
$ws.Cells.Item(122, 13).Value = -6166.999899999999

# Row 126
$ws.Cells.Item(126, 8).Value = 5184.5938
$ws.Cells.Item(126, 9).Value = 3995.9565
$ws.Cells.Item(126, 10).Value = 8222.223
$ws.Cells.Item(126, 11).Value = 11987.8695
$ws.Cells.Item(126, 12).Value = 24666.669
$ws.Cells.Item(126, 13).Value = -9517.869499999999
$ws.Cells.Item(126, 14).Value = -29606.669

# Row 132
$ws.Cells.Item(132, 8).Value = 5130.114
$ws.Cells.Item(132, 9).Value = 3269.5
$ws.Cells.Item(132, 10).Value = 6990.727
$ws.Cells.Item(132, 11).Value = 9808.5
$ws.Cells.Item(132, 12).Value = 20972.181
$ws.Cells.Item(132, 13).Value = -7278.5
$ws.Cells.Item(132, 14).Value = -26032.181

$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Cells.Item(132, 8).Value = 10804.462
$ws.Cells.Item(132, 9).Value = 100000
$ws.Cells.Item(132, 10).Value = 3371.5
$ws.Cells.Item(132, 11).Value = 300000
$ws.Cells.Item(132, 12).Value = 10114.5
$ws.Cells.Item(132, 13).Value = -297470
$ws.Cells.Item(132, 14).Value = -15174.5
